## Add a new "vaccination_transfer" worksheet (copy of "vaccination" with an
## extra "inAppDeliveryCode" column) after "vaccination_print".

$wb = $excel.ActiveWorkbook

# Create the new worksheet as the last tab, then rename it.
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "vaccination_transfer"

# Resolve the source sheet reference AFTER adding the new sheet (adding a
# sheet can invalidate previously-held worksheet references in this
# runtime).
$src = $wb.Worksheets.Item("vaccination")

# Copy the two data rows (header + one data row) across columns A:I.
for ($c = 1; $c -le 9; $c++) {
    $ws.Cells.Item(1, $c).Value = $src.Cells.Item(1, $c).Value2
    $ws.Cells.Item(2, $c).Value = $src.Cells.Item(2, $c).Value2
}

# Re-apply the date number format on columns C and H (same two columns that
# carry it on the source sheet) so the engine reuses the existing style
# instead of minting a new one.
$ws.Cells.Item(1, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(2, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(1, 8).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(2, 8).NumberFormat = "yyyy\-mm\-dd"

# New column J: inAppDeliveryCode.
$ws.Cells.Item(1, 10).Value = "inAppDeliveryCode"
$ws.Cells.Item(2, 10).Value = "Y8P8ECFN8"
$ws.Cells.Item(1, 10).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(2, 10).NumberFormat = "yyyy\-mm\-dd"

# Match column widths to the source sheet for columns A:I ...
for ($c = 1; $c -le 9; $c++) {
    $ws.Columns.Item($c).ColumnWidth = $src.Columns.Item($c).ColumnWidth
}
# ... and give column J the same width as the similarly-sized
# "cantonCodeSender" column on vaccination_print.
$wsPrint = $wb.Worksheets.Item("vaccination_print")
$ws.Columns.Item(10).ColumnWidth = $wsPrint.Columns.Item(8).ColumnWidth

# Make the new sheet the active tab / selection, matching the authored file.
$ws.Activate()
$ws.Range("I3").Select() | Out-Null
